# Edit script: reorganize rows 2-22 in Sheet1 of the workbook.
# The diff shows that rows 2-22 are rotated (the last 4 rows move to the
# front) and the 4 rows that moved (EVOL3088, EVOL3087, EVOL3086, EVOL3089)
# receive brand new Google Drive ids/links, while the rest keep their
# existing ids/links, just shifted to new row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2, 'EVOL3088.png', '1eIiLw76jRpOSREMEBPIDpjFFqkxjchM2', 'https://drive.google.com/file/d/1eIiLw76jRpOSREMEBPIDpjFFqkxjchM2/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1eIiLw76jRpOSREMEBPIDpjFFqkxjchM2', 'EVOL3088'),
    @(3, 'EVOL3087.png', '1po1vu4rn4nYzrCkc62BP4VTXYT0ECK8h', 'https://drive.google.com/file/d/1po1vu4rn4nYzrCkc62BP4VTXYT0ECK8h/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1po1vu4rn4nYzrCkc62BP4VTXYT0ECK8h', 'EVOL3087'),
    @(4, 'EVOL3086.png', '1BvURk9_AZvhv7bu5cU165dzxtk54lrfP', 'https://drive.google.com/file/d/1BvURk9_AZvhv7bu5cU165dzxtk54lrfP/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1BvURk9_AZvhv7bu5cU165dzxtk54lrfP', 'EVOL3086'),
    @(5, 'EVOL3089.png', '1lkFpCsC8W2m3LUIiKnNyq5iQ5T4dSGWY', 'https://drive.google.com/file/d/1lkFpCsC8W2m3LUIiKnNyq5iQ5T4dSGWY/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1lkFpCsC8W2m3LUIiKnNyq5iQ5T4dSGWY', 'EVOL3089'),
    @(6, 'EVOL0011.png', '1SQ4RO0DwgN6iB16qKKJW4nOKO-P-9AZ3', 'https://drive.google.com/file/d/1SQ4RO0DwgN6iB16qKKJW4nOKO-P-9AZ3/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1SQ4RO0DwgN6iB16qKKJW4nOKO-P-9AZ3', 'EVOL0011'),
    @(7, 'EVOL0010.png', '1GxnmEGjpuhLnIQTUrBPOHbaKsMVc8ztN', 'https://drive.google.com/file/d/1GxnmEGjpuhLnIQTUrBPOHbaKsMVc8ztN/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1GxnmEGjpuhLnIQTUrBPOHbaKsMVc8ztN', 'EVOL0010'),
    @(8, 'EVOL0009.png', '17ABzRsKQ-aAb7OOIy4gjt7EM2XOgeKi3', 'https://drive.google.com/file/d/17ABzRsKQ-aAb7OOIy4gjt7EM2XOgeKi3/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=17ABzRsKQ-aAb7OOIy4gjt7EM2XOgeKi3', 'EVOL0009'),
    @(9, 'EVOL0008.png', '1io205s0lr5DzUCVGViGbT2u_eoQJZb2C', 'https://drive.google.com/file/d/1io205s0lr5DzUCVGViGbT2u_eoQJZb2C/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1io205s0lr5DzUCVGViGbT2u_eoQJZb2C', 'EVOL0008'),
    @(10, 'EVOL6765.jpg', '15gfbTSnSVxp3RepFCtm2oOGlz219oyZh', 'https://drive.google.com/file/d/15gfbTSnSVxp3RepFCtm2oOGlz219oyZh/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=15gfbTSnSVxp3RepFCtm2oOGlz219oyZh', 'EVOL6765'),
    @(11, 'EVOL6760.png', '1vn8sI2ITt5XjGHLb23Z05vILv_DVLw05', 'https://drive.google.com/file/d/1vn8sI2ITt5XjGHLb23Z05vILv_DVLw05/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1vn8sI2ITt5XjGHLb23Z05vILv_DVLw05', 'EVOL6760'),
    @(12, 'EVOL6715.png', '1SN4-JQEpjlD5Q-TDn2v5lNQ1hFyKVot1', 'https://drive.google.com/file/d/1SN4-JQEpjlD5Q-TDn2v5lNQ1hFyKVot1/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1SN4-JQEpjlD5Q-TDn2v5lNQ1hFyKVot1', 'EVOL6715'),
    @(13, 'EVOL6222.png', '1O-xOhjMBTW4WINz5wsnhNFHCn5EnT0Ki', 'https://drive.google.com/file/d/1O-xOhjMBTW4WINz5wsnhNFHCn5EnT0Ki/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1O-xOhjMBTW4WINz5wsnhNFHCn5EnT0Ki', 'EVOL6222'),
    @(14, 'EVOL6221.png', '19iuIRdeE_RammaNztz1KygdNq1tGTL8C', 'https://drive.google.com/file/d/19iuIRdeE_RammaNztz1KygdNq1tGTL8C/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=19iuIRdeE_RammaNztz1KygdNq1tGTL8C', 'EVOL6221'),
    @(15, 'EVOL6210.png', '15mNv1wWkmsHpiLiWlg8jgv5yQVZt4Z1l', 'https://drive.google.com/file/d/15mNv1wWkmsHpiLiWlg8jgv5yQVZt4Z1l/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=15mNv1wWkmsHpiLiWlg8jgv5yQVZt4Z1l', 'EVOL6210'),
    @(16, 'EVOL6205.png', '11hBmjH-q6kh0_CNOGWTw6LnsXjbves0f', 'https://drive.google.com/file/d/11hBmjH-q6kh0_CNOGWTw6LnsXjbves0f/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=11hBmjH-q6kh0_CNOGWTw6LnsXjbves0f', 'EVOL6205'),
    @(17, 'EVOL3957.png', '14c8KBXNh6nrKA5S4u67BmU7nOAZL4Adk', 'https://drive.google.com/file/d/14c8KBXNh6nrKA5S4u67BmU7nOAZL4Adk/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=14c8KBXNh6nrKA5S4u67BmU7nOAZL4Adk', 'EVOL3957'),
    @(18, 'EVOL3953.png', '1PWg5nUjF3Lzsf6MSF9eUD0K364TAsO1O', 'https://drive.google.com/file/d/1PWg5nUjF3Lzsf6MSF9eUD0K364TAsO1O/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1PWg5nUjF3Lzsf6MSF9eUD0K364TAsO1O', 'EVOL3953'),
    @(19, 'EVOL3959.png', '1MtgmhF4llFOz-8-Kg3Xk0eOPlIrE6NqZ', 'https://drive.google.com/file/d/1MtgmhF4llFOz-8-Kg3Xk0eOPlIrE6NqZ/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1MtgmhF4llFOz-8-Kg3Xk0eOPlIrE6NqZ', 'EVOL3959'),
    @(20, 'EVOL3955.png', '1sRu-iyXG_A4UpHTtmWy4m_O-1eVKvO-f', 'https://drive.google.com/file/d/1sRu-iyXG_A4UpHTtmWy4m_O-1eVKvO-f/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1sRu-iyXG_A4UpHTtmWy4m_O-1eVKvO-f', 'EVOL3955'),
    @(21, 'EVOL3961.png', '1h2_X5ZaWr3LDrgj_YErdS25QazR7X0od', 'https://drive.google.com/file/d/1h2_X5ZaWr3LDrgj_YErdS25QazR7X0od/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1h2_X5ZaWr3LDrgj_YErdS25QazR7X0od', 'EVOL3961'),
    @(22, 'EVOL0440.png', '1Ms54HbfpRLBZCKT8PoD-2ajd7x6LXc-z', 'https://drive.google.com/file/d/1Ms54HbfpRLBZCKT8PoD-2ajd7x6LXc-z/view?usp=drivesdk', 'https://drive.google.com/uc?export=view&id=1Ms54HbfpRLBZCKT8PoD-2ajd7x6LXc-z', 'EVOL0440')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
